# Fruta / hortaliza, semanal
# Insert one new daily price record as a new row 54, pushing every
# existing record (old rows 54..157) down by one row (new rows 55..158).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 54 downward by inserting a fresh blank row at position 54.
$ws.Rows(54).Insert()

# Populate the newly inserted row 54 with the new Cereza price record.
$ws.Cells.Item(54, 1).Value = 7
$ws.Cells.Item(54, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(54, 3).Value = "Ñuble"
$ws.Cells.Item(54, 4).Value = 44935
$ws.Cells.Item(54, 5).Value = 16
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100103
$ws.Cells.Item(54, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(54, 9).Value = 100103001
$ws.Cells.Item(54, 10).Value = "Cereza"
$ws.Cells.Item(54, 11).Value = "Lapins"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 60
$ws.Cells.Item(54, 14).Value = 5000
$ws.Cells.Item(54, 15).Value = 5500
$ws.Cells.Item(54, 16).Value = 5250
$ws.Cells.Item(54, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(54, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(54, 19).Value = 525
$ws.Cells.Item(54, 20).Value = 10
